# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps recorded in the zh-cn and de-de handback status sheets.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-14 04:13:51"
$wsZh.Range("E5").Value = "2016-03-14 04:13:51"
$wsZh.Range("H3").Value = "2016-03-14 04:14:10"
$wsZh.Range("H5").Value = "2016-03-14 04:14:10"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-14 04:13:54"
$wsDe.Range("E5").Value = "2016-03-14 04:13:54"
$wsDe.Range("H3").Value = "2016-03-14 04:14:15"
$wsDe.Range("H5").Value = "2016-03-14 04:14:15"
